# Applies odds/clean-sheet value updates to rows 3-5 per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("H3").Value = 4  # was 3.9
$ws.Range("J3").Value = 2.2  # was 2.25
$ws.Range("K3").Value = 2.5  # was 2.4
$ws.Range("L3").Value = 4.33  # was 4.5
$ws.Range("N3").Value = 17  # was 15
$ws.Range("O3").Value = 1.14  # was 1.18
$ws.Range("P3").Value = 5.5  # was 4.5
$ws.Range("Q3").Value = 1.5  # was 1.62
$ws.Range("R3").Value = 2.5  # was 2.25
$ws.Range("S3").Value = 1.25  # was 1.29
$ws.Range("T3").Value = 3.75  # was 3.5
$ws.Range("U3").Value = 1.5  # was 1.57
$ws.Range("V3").Value = 2.5  # was 2.25
$ws.Range("W3").Value = 11  # was 9.5
$ws.Range("X3").Value = 10  # was 9.5
$ws.Range("AA3").Value = 12  # was 13
$ws.Range("AB3").Value = 19  # was 21
$ws.Range("AC3").Value = 19  # was 15
$ws.Range("AE3").Value = 12  # was 13
$ws.Range("AF3").Value = 34  # was 41
$ws.Range("AG3").Value = 101  # was 126
$ws.Range("AH3").Value = 19  # was 17
$ws.Range("AI3").Value = 29  # was 26
$ws.Range("AK3").Value = 51  # was 41
$ws.Range("AL3").Value = 29  # was 34
$ws.Range("AM3").Value = 29  # was 34
$ws.Range("AP3").Value = 15  # was 17
$ws.Range("AQ3").Value = 23  # was 26
$ws.Range("AS3").Value = 81  # was 101
$ws.Range("AT3").Value = 3.75  # was 3.5
$ws.Range("AY3").Value = 23  # was 26
$ws.Range("BA3").Value = 67  # was 81
$ws.Range("BB3").Value = 126  # was 151

# Row 4
$ws.Range("G4").Value = 2.2  # was 2.1
$ws.Range("I4").Value = 3.3  # was 3.7
$ws.Range("J4").Value = 2.88  # was 2.75
$ws.Range("L4").Value = 3.75  # was 4
$ws.Range("M4").Value = 1.05  # was 1.07
$ws.Range("N4").Value = 11  # was 9
$ws.Range("O4").Value = 1.25  # was 1.33
$ws.Range("P4").Value = 3.75  # was 3.25
$ws.Range("Q4").Value = 1.85  # was 2.05
$ws.Range("R4").Value = 1.95  # was 1.75
$ws.Range("S4").Value = 1.4  # was 1.44
$ws.Range("T4").Value = 2.75  # was 2.63
$ws.Range("U4").Value = 1.73  # was 1.83
$ws.Range("V4").Value = 2  # was 1.83
$ws.Range("W4").Value = 8  # was 7
$ws.Range("X4").Value = 11  # was 9.5
$ws.Range("Z4").Value = 21  # was 19
$ws.Range("AB4").Value = 26  # was 29
$ws.Range("AC4").Value = 10  # was 9
$ws.Range("AE4").Value = 13  # was 15
$ws.Range("AF4").Value = 41  # was 51
$ws.Range("AG4").Value = 201  # was 251
$ws.Range("AH4").Value = 11  # was 10
$ws.Range("AI4").Value = 17  # was 19
$ws.Range("AJ4").Value = 12  # was 13
$ws.Range("AK4").Value = 34  # was 41
$ws.Range("AL4").Value = 26  # was 29
$ws.Range("AM4").Value = 34  # was 41
$ws.Range("AN4").Value = 4.33  # was 4
$ws.Range("AP4").Value = 21  # was 23
$ws.Range("AT4").Value = 2.75  # was 2.63
$ws.Range("AW4").Value = 5  # was 5.5
$ws.Range("AX4").Value = 17  # was 21
$ws.Range("AY4").Value = 26  # was 29
$ws.Range("AZ4").Value = 51  # was 67
$ws.Range("BA4").Value = 81  # was 101
$ws.Range("BB4").Value = 151  # was 201

# Row 5
$ws.Range("G5").Value = 1.91  # was 2.1
$ws.Range("H5").Value = 3.4  # was 3.25
$ws.Range("I5").Value = 4.1  # was 3.5
$ws.Range("J5").Value = 2.5  # was 2.75
$ws.Range("K5").Value = 2.3  # was 2.25
$ws.Range("L5").Value = 4  # was 3.75
$ws.Range("O5").Value = 1.18  # was 1.22
$ws.Range("P5").Value = 4.5  # was 4
$ws.Range("Q5").Value = 1.65  # was 1.75
$ws.Range("R5").Value = 2.2  # was 2.05
$ws.Range("S5").Value = 1.3  # was 1.33
$ws.Range("T5").Value = 3.4  # was 3.25
$ws.Range("Z5").Value = 17  # was 19
$ws.Range("AA5").Value = 13  # was 15
$ws.Range("AB5").Value = 21  # was 23
$ws.Range("AC5").Value = 13  # was 12
$ws.Range("AD5").Value = 7  # was 6.5
$ws.Range("AE5").Value = 12  # was 11
$ws.Range("AH5").Value = 15  # was 13
$ws.Range("AI5").Value = 23  # was 19
$ws.Range("AL5").Value = 29  # was 26
$ws.Range("AN5").Value = 4  # was 4.33
$ws.Range("AO5").Value = 10  # was 11
$ws.Range("AP5").Value = 17  # was 19
$ws.Range("AR5").Value = 41  # was 51
$ws.Range("AS5").Value = 101  # was 126
$ws.Range("AT5").Value = 3.4  # was 3.25
$ws.Range("AV5").Value = 41  # was 51
$ws.Range("AW5").Value = 6  # was 5.5
$ws.Range("AX5").Value = 19  # was 17
